$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3065.5
$ws.Range("I32").Value = 1225
$ws.Range("J32").Value = 3591.3572
$ws.Range("K32").Value = 1225
$ws.Range("L32").Value = 3591.3572
$ws.Range("M32").Value = -899
$ws.Range("N32").Value = -4243.3572

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 76927070
$ws.Range("J113").Value = 4440.6665
$ws.Range("L113").Value = 4440.6665
$ws.Range("N113").Value = -10948.6665

# Sheet ALC, row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 12530
$ws.Range("I131").Value = 7000
$ws.Range("J131").Value = 13144.444
$ws.Range("K131").Value = 21000
$ws.Range("L131").Value = 39433.33199999999
$ws.Range("M131").Value = -15960
$ws.Range("N131").Value = -49513.33199999999

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1398.7368
$ws.Range("I2").Value = 1316.2354
$ws.Range("K2").Value = 1316.2354
$ws.Range("M2").Value = -1203.2354

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1318.5139
$ws.Range("I32").Value = 1266.662
$ws.Range("K32").Value = 1266.662
$ws.Range("M32").Value = -979.662

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15875372
$ws.Range("I61").Value = 23811150
$ws.Range("K61").Value = 23811150
$ws.Range("M61").Value = -23810938

# Sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 66670570
$ws.Range("I63").Value = 166667980
$ws.Range("J63").Value = 16671858
$ws.Range("K63").Value = 166667980
$ws.Range("L63").Value = 16671858
$ws.Range("M63").Value = -166667294
$ws.Range("N63").Value = -16673230

# Sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 66670570
$ws.Range("I66").Value = 166667980
$ws.Range("J66").Value = 16671858
$ws.Range("K66").Value = 833339900
$ws.Range("L66").Value = 83359290
$ws.Range("M66").Value = -833336468
$ws.Range("N66").Value = -83366154

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3539.4443
$ws.Range("I74").Value = 3660.5
$ws.Range("J74").Value = 3297.3333
$ws.Range("K74").Value = 3660.5
$ws.Range("L74").Value = 3297.3333
$ws.Range("M74").Value = -2786.5
$ws.Range("N74").Value = -5045.3333

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3539.4443
$ws.Range("I77").Value = 3660.5
$ws.Range("J77").Value = 3297.3333
$ws.Range("K77").Value = 18302.5
$ws.Range("L77").Value = 16486.6665
$ws.Range("M77").Value = -13934.5
$ws.Range("N77").Value = -25222.6665

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6945881
$ws.Range("I88").Value = 13889450
$ws.Range("J88").Value = 2311.9167
$ws.Range("K88").Value = 13889450
$ws.Range("L88").Value = 2311.9167
$ws.Range("M88").Value = -13889044
$ws.Range("N88").Value = -3123.9167

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 6945881
$ws.Range("I91").Value = 13889450
$ws.Range("J91").Value = 2311.9167
$ws.Range("K91").Value = 13889450
$ws.Range("L91").Value = 2311.9167
$ws.Range("M91").Value = -13888046
$ws.Range("N91").Value = -5119.9167

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1398.7368
$ws.Range("I116").Value = 1316.2354
$ws.Range("K116").Value = 1316.2354
$ws.Range("M116").Value = 977.7646

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 15876385
$ws.Range("I122").Value = 22224938
$ws.Range("K122").Value = 66674814
$ws.Range("M122").Value = -66672364

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 38463364
$ws.Range("I132").Value = 41668480
$ws.Range("K132").Value = 125005440
$ws.Range("M132").Value = -125002910

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 15875372
$ws.Range("I136").Value = 23811150
$ws.Range("K136").Value = 71433450
$ws.Range("M136").Value = -71430900

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1398.7368
$ws.Range("I3").Value = 1316.2354
$ws.Range("K3").Value = 1316.2354
$ws.Range("M3").Value = -1202.2354

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11681.363
$ws.Range("I20").Value = 16850.666
$ws.Range("J20").Value = 5478.2
$ws.Range("K20").Value = 16850.666
$ws.Range("L20").Value = 5478.2
$ws.Range("M20").Value = -16603.666
$ws.Range("N20").Value = -5972.2

# Sheet BSM, row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35586.688
$ws.Range("I82").Value = 8431.25
$ws.Range("J82").Value = 62742.125
$ws.Range("K82").Value = 8431.25
$ws.Range("L82").Value = 62742.125
$ws.Range("M82").Value = -8048.25
$ws.Range("N82").Value = -63508.125

# Sheet BSM, row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 35586.688
$ws.Range("I85").Value = 8431.25
$ws.Range("J85").Value = 62742.125
$ws.Range("K85").Value = 8431.25
$ws.Range("L85").Value = 62742.125
$ws.Range("M85").Value = -7105.25
$ws.Range("N85").Value = -65394.125

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3106.423
$ws.Range("I86").Value = 2921.4443
$ws.Range("K86").Value = 2921.4443
$ws.Range("M86").Value = -1798.4443

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3106.423
$ws.Range("I89").Value = 2921.4443
$ws.Range("K89").Value = 2921.4443
$ws.Range("M89").Value = -8991.2215

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 746.2308
$ws.Range("I99").Value = 626.4545000000001
$ws.Range("J99").Value = 1405
$ws.Range("K99").Value = 626.4545000000001
$ws.Range("L99").Value = 1405
$ws.Range("M99").Value = 871.5454999999999
$ws.Range("N99").Value = -4401

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2253.6
$ws.Range("J132").Value = 2756.1667
$ws.Range("L132").Value = 24805.5003
$ws.Range("N132").Value = -29865.5003

# Sheet GSM, row 49
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 40000
$ws.Range("J49").Value = 40000
$ws.Range("L49").Value = 40000
$ws.Range("N49").Value = -40368

# Sheet GSM, row 92
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1783.0526
$ws.Range("I97").Value = 1717.3914
$ws.Range("K97").Value = 1717.3914
$ws.Range("M97").Value = -1221.3914

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 31252456
$ws.Range("I122").Value = 1755.1818
$ws.Range("K122").Value = 5265.5454
$ws.Range("M122").Value = -2815.5454

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1815.2
$ws.Range("I40").Value = 1109
$ws.Range("K40").Value = 1109
$ws.Range("M40").Value = -973

# Sheet LTW, row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4735586.5
$ws.Range("I43").Value = 3156623.5
$ws.Range("K43").Value = 3156623.5
$ws.Range("M43").Value = -3156430.5

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6649.2104
$ws.Range("I61").Value = 4343.8
$ws.Range("J61").Value = 15294.5
$ws.Range("K61").Value = 4343.8
$ws.Range("L61").Value = 15294.5
$ws.Range("M61").Value = -4141.8
$ws.Range("N61").Value = -15698.5

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3182.9375
$ws.Range("I82").Value = 3495.2307
$ws.Range("K82").Value = 3495.2307
$ws.Range("M82").Value = -3134.2307

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3182.9375
$ws.Range("I85").Value = 3495.2307
$ws.Range("K85").Value = 3495.2307
$ws.Range("M85").Value = -2247.2307

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6649.2104
$ws.Range("I113").Value = 4343.8
$ws.Range("J113").Value = 15294.5
$ws.Range("K113").Value = 4343.8
$ws.Range("L113").Value = 15294.5
$ws.Range("M113").Value = -2173.8
$ws.Range("N113").Value = -19634.5

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3336.9355
$ws.Range("I122").Value = 2633.8635
$ws.Range("K122").Value = 7901.5905
$ws.Range("M122").Value = -5451.5905

# Sheet WVR, row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2763.6924
$ws.Range("I23").Value = 2196.4443
$ws.Range("J23").Value = 4040
$ws.Range("K23").Value = 2196.4443
$ws.Range("L23").Value = 4040
$ws.Range("M23").Value = -1967.4443
$ws.Range("N23").Value = -4498

# Sheet WVR, row 55
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 11184.333
$ws.Range("I55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("M55").Value = -9723

# Sheet WVR, row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 35249.75
$ws.Range("J104").Value = 35249.75
$ws.Range("L104").Value = 35249.75
$ws.Range("N104").Value = -42237.75

# Sheet WVR, row 110
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Sheet WVR, row 111
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
